# Applies the cryptos.xlsx price/volume refresh described in the commit:
# "Updated cryptos list on Wed Mar 27 14:25:45 UTC 2024 with GitHub Actions"
#
# All Coin/Link/Price/Volume(1h) cells are stored as plain text in the workbook,
# so numeric-looking Price values are written with a leading apostrophe to force
# Excel to keep them as text instead of auto-converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.503.28'
$ws.Range('E2').Value = '  -0.99%  '
$ws.Range('D3').Value = '3.530.91'
$ws.Range('E3').Value = '  -1.68%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('D5').Value = "'570.36"
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('D6').Value = "'184.00"
$ws.Range('E6').Value = '  -4.38%  '
$ws.Range('D7').Value = '3.529.44'
$ws.Range('E7').Value = '  -1.59%  '
$ws.Range('D8').Value = "'0.614"
$ws.Range('E8').Value = '  -3.61%  '
$ws.Range('E9').Value = '  -0.10%  '
$ws.Range('D10').Value = "'0.183"
$ws.Range('E10').Value = '  +1.66%  '
$ws.Range('D11').Value = "'0.646"
$ws.Range('E11').Value = '  -2.73%  '
$ws.Range('D12').Value = "'53.92"
$ws.Range('E12').Value = '  -5.50%  '
$ws.Range('D13').Value = "'0.0000299"
$ws.Range('E13').Value = '  -2.00%  '
$ws.Range('D14').Value = "'9.48"
$ws.Range('E14').Value = '  -3.11%  '
$ws.Range('D15').Value = '4.096.29'
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('D16').Value = "'19.34"
$ws.Range('E16').Value = '  -4.44%  '
$ws.Range('D17').Value = '3.517.15'
$ws.Range('E17').Value = '  -2.50%  '
$ws.Range('D18').Value = '69.175.53'
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').Value = "'12.32"
$ws.Range('E19').Value = '  -1.98%  '
$ws.Range('E20').Value = '  -1.68%  '
$ws.Range('D21').Value = "'1.04"
$ws.Range('E21').Value = '  -1.35%  '
$ws.Range('D22').Value = "'505.84"
$ws.Range('E22').Value = '  +4.53%  '
$ws.Range('D23').Value = "'19.50"
$ws.Range('E23').Value = '  -0.97%  '
$ws.Range('D24').Value = "'4.91"
$ws.Range('E24').Value = '  -3.93%  '
$ws.Range('D25').Value = "'4.33"
$ws.Range('E25').Value = '  -1.69%  '
$ws.Range('D26').Value = "'94.06"
$ws.Range('E26').Value = '  +4.67%  '
$ws.Range('D27').Value = "'11.32"
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('D28').Value = "'2.93"
$ws.Range('E28').Value = '  -5.64%  '
$ws.Range('D29').Value = "'9.18"
$ws.Range('E29').Value = '  -2.56%  '
$ws.Range('D30').Value = "'31.45"
$ws.Range('E30').Value = '  -2.98%  '
$ws.Range('D31').Value = "'7.55"
$ws.Range('E31').Value = '  -3.33%  '
$ws.Range('D32').Value = "'12.42"
$ws.Range('E32').Value = '  +1.98%  '
$ws.Range('D33').Value = "'65.20"
$ws.Range('E33').Value = '  -1.59%  '
$ws.Range('D34').Value = "'0.115"
$ws.Range('E34').Value = '  -5.93%  '
$ws.Range('D35').Value = "'568.82"
$ws.Range('E35').Value = '  -6.47%  '
$ws.Range('D36').Value = "'3.14"
$ws.Range('E36').Value = '  +4.98%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D37').Value = "'37.96"
$ws.Range('E37').Value = '  -5.51%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').Value = "'1.00"
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('D39').Value = "'0.402"
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0777'
$ws.Range('E40').Value = '  -5.14%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = "'3.20"
$ws.Range('E41').Value = '  +0.90%  '
$ws.Range('D42').Value = "'3.39"
$ws.Range('E42').Value = '  -4.42%  '
$ws.Range('E43').Value = '  -9.53%  '
$ws.Range('D44').Value = "'3.54"
$ws.Range('E44').Value = '  +4.48%  '
$ws.Range('D45').Value = "'2.98"
$ws.Range('E45').Value = '  -5.42%  '
$ws.Range('D46').Value = "'0.0445"
$ws.Range('E46').Value = '  -1.82%  '
$ws.Range('D47').Value = '3.170.01'
$ws.Range('E47').Value = '  -4.35%  '
$ws.Range('D48').Value = "'9.30"
$ws.Range('E48').Value = '  -3.07%  '
$ws.Range('D49').Value = "'0.135"
$ws.Range('E49').Value = '  -2.47%  '
$ws.Range('B50').Value = 'LidoDAOToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D50').Value = "'3.25"
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('B51').Value = 'FirstDigitalUSD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D51').Value = "'0.996"
$ws.Range('E51').Value = '  -0.40%  '
